# Update cryptocurrency price (column D) and volume change (column E) values
# as scraped on Mon Aug 21 09:54:38 UTC 2023 with GitHub Actions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> [Price, Volume(1h)] new values.
# Only rows whose Price or Volume text actually changed are listed with new
# values; rows whose Price stayed the same keep their original Price text.
$updates = @{
    2  = @("26.187.51", "  -0.44%  ")
    3  = @("1.679.87",  "  -0.02%  ")
    4  = @("1.006",     "  -0.03%  ")
    5  = @("215.74",    "  -1.05%  ")
    6  = @("0.5257",    "  -2.58%  ")
    7  = @("1.006",     "  -0.07%  ")
    8  = @("0.2684",    "  -0.17%  ")
    9  = @("0.06357",   "  -1.84%  ")
    10 = @("21.40",     "  -2.65%  ")
    11 = @("0.07617",   "  +1.04%  ")
    12 = @("1.687.54",  "  +0.38%  ")
    13 = @("4.523",     "  -0.15%  ")
    14 = @("0.5725",    "  -1.09%  ")
    15 = @("0.000008203","  -2.93%  ")
    16 = @("66.34",     "  +2.36%  ")
    17 = @("26.234.67", "  -0.32%  ")
    18 = @("1.006",     "  -0.03%  ")
    19 = @("4.859",     "  -1.00%  ")
    20 = @("10.72",     "  -1.54%  ")
    21 = @("189.75",    "  -0.65%  ")
    22 = @("6.219",     "  +0.02%  ")
    23 = @("1.007",     "  -0.05%  ")
    24 = @("149.07",    "  +1.95%  ")
    25 = @("0.1256",    "  -2.77%  ")
    26 = @("7.700",     "  -1.76%  ")
    27 = @("15.87",     "  +0.71%  ")
    28 = @("0.06419",   "  -1.24%  ")
    29 = @("1.371",     "  -1.54%  ")
    30 = @("1.312",     "  -0.75%  ")
    31 = @("3.564",     "  -0.46%  ")
    32 = @("3.556",     "  -0.54%  ")
    33 = @("1.674",     "  +0.48%  ")
    34 = @("1.014",     "  -2.02%  ")
    35 = @("0.6081",    "  -1.31%  ")
    36 = @("2.419",     "  +0.76%  ")
    37 = @("2.743",     "  +0.78%  ")
    38 = @("0.01641",   "  +1.35%  ")
    39 = @("6.160",     "  -1.38%  ")
    40 = @("1.092.40",  "  -1.75%  ")
    41 = @("0.8808",    "  +1.17%  ")
    42 = @("1.010",     "  -0.51%  ")
    43 = @("100.42",    "  -0.12%  ")
    44 = @("1.833.46",  "  +0.23%  ")
    45 = @("57.40",     "  +0.37%  ")
    46 = @("0.00000000108","  -1.13%  ")
    47 = @("0.9988",    "  -0.12%  ")
    48 = @("8.042",     "  -1.20%  ")
    49 = @("0.05264",   "  -0.17%  ")
    50 = @("0.4281",    "  -0.18%  ")
    51 = @("5.996",     "  -1.32%  ")
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]

    $dCell = $ws.Cells.Item($row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $values[0]

    $eCell = $ws.Cells.Item($row, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $values[1]
}
